$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "6b38ddb6-1cdc-5f7b-6efa-d9c911cf4972"
$ws.Range("B3").Value = "375e8f77-3102-4813-5cc4-b9aac0c0d908"
$ws.Range("E2").Value = "5c1a8653-7790-743f-c400-73728ad2d17f"
$ws.Range("E3").Value = "724cda7e-3e73-6365-b639-7fa89fdac40f"

$ws.Columns("E").ColumnWidth = 36.5

$ws.Range("G2").Select()
